$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "Niccoló Orsi"
$ws.Range("B29").Value = "Stefano Tita | Clitoriders"
$ws.Range("C29").Value = "Federico  Manica | iMontagna"
$ws.Range("D29").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E29").Value = "Mattia Festi | Shark Attack"
$ws.Range("F29").Value = "Michael Bertè | A.C. Denti"
